$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERTS")

# Columns B-F are plain text / time strings that Excel's type-inference leaves
# as text already, so a direct .Value assignment is safe for them.
#
# Column A holds ISO-style "YYYY-MM-DD" strings. A direct .Value assignment on
# such a string gets auto-parsed into a date serial by Excel's smart-entry
# heuristic, which would also change the cell's number format / style. To
# preserve the plain-text representation used by the existing rows (and avoid
# minting a new, unused style in styles.xml), stage the literal text as a
# formula result in a scratch cell, then copy/paste-special (values-only) it
# into the target cell - paste-special values bypass the smart re-parsing.

$scratch = $ws.Range("Z1")

function Set-TextDate($cellAddr, $text) {
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextDate "A6" "2026-01-31"
$ws.Range("B6").Value = "22:35:01"
$ws.Range("C6").Value = "22:00"
$ws.Range("D6").Value = "Living Room"
$ws.Range("E6").Value = "CRITICAL"
$ws.Range("F6").Value = "FALL_DETECTED"

Set-TextDate "A7" "2026-01-31"
$ws.Range("B7").Value = "22:35:04"
$ws.Range("C7").Value = "22:00"
$ws.Range("D7").Value = "Living Room"
$ws.Range("E7").Value = "CRITICAL"
$ws.Range("F7").Value = "FALL_DETECTED"

$scratch.Clear()
$excel.CutCopyMode = $false
